# Import OBI term 'transcription factor binding site identification'
#
# The OBI_input worksheet lists ontology terms ordered by their OBI
# numeric identifier. A new term (OBI_0000291) is inserted between the
# existing OBI_0000281 row (row 31) and the OBI_0000341 row (old row 32),
# pushing every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a fresh row at position 32; everything below shifts down by one.
$ws.Rows.Item(32).Insert()

# Populate the new row's data (columns: A = source IRI, B = source label,
# C = "Include in View" flag; no preferred-label override in column D).
$ws.Range("A32").Value = "http://purl.obolibrary.org/obo/OBI_0000291"
$ws.Range("B32").Value = "transcription factor binding site  identification"
$ws.Range("C32").Value = "y"

# Match the style used by the surrounding rows in column C.
$ws.Range("C32").Font.Name = $ws.Range("C33").Font.Name
$ws.Range("C32").Font.Size = $ws.Range("C33").Font.Size

# Leave the view scrolled/selected roughly where the author left it.
$ws.Range("D29").Select()
